# Commit: "added clutter counterbalance sheet, clutter change feature"
#
# Adds a new clutterChangeEnabled flag plus four generic counter columns
# (a, b, c, d) in L1:P1 of Sheet1, and re-derives the per-row values:
#   - rows 2-5 ("training"/"lexical" rows): word/nonword trial counts
#     (E/F) are collapsed to 1 and their previous values are kept in the
#     new c/d columns (O/P); clutterChangeEnabled = FALSE, a/b = 1.
#   - rows 6-9 ("full task" rows): iti_min/iti_max (C/D) are collapsed to
#     1 and their previous values are kept in the new a/b columns (M/N);
#     clutterChangeEnabled = TRUE. Row 6 also updates its own E/F from
#     8 to 20 (and that new value is mirrored into c/d).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header row (L1:P1) -----------------------------------------------
$ws.Range("L1").Value = "clutterChangeEnabled"
$ws.Range("M1").Value = "a"
$ws.Range("N1").Value = "b"
$ws.Range("O1").Value = "c"
$ws.Range("P1").Value = "d"

# --- Row 2 (training_lexical / Georgia) -----------------------------------
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 1
$ws.Range("L2").Value = $false
$ws.Range("M2").Value = 1
$ws.Range("N2").Value = 1
$ws.Range("O2").Value = 20
$ws.Range("P2").Value = 20

# --- Row 3 (lexical_wo_driving_roboto) -------------------------------------
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 1
$ws.Range("L3").Value = $false
$ws.Range("M3").Value = 1
$ws.Range("N3").Value = 1
$ws.Range("O3").Value = 30
$ws.Range("P3").Value = 30

# --- Row 4 (lexical_wo_driving_neuefrutigerworld) --------------------------
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 1
$ws.Range("L4").Value = $false
$ws.Range("M4").Value = 1
$ws.Range("N4").Value = 1
$ws.Range("O4").Value = 30
$ws.Range("P4").Value = 30

# --- Row 5 (lexical_wo_driving_eurostile) ----------------------------------
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 1
$ws.Range("L5").Value = $false
$ws.Range("M5").Value = 1
$ws.Range("N5").Value = 1
$ws.Range("O5").Value = 30
$ws.Range("P5").Value = 30

# --- Row 6 (full_task_training / Georgia) ----------------------------------
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 20
$ws.Range("F6").Value = 20
$ws.Range("L6").Value = $true
$ws.Range("M6").Value = 5
$ws.Range("N6").Value = 10
$ws.Range("O6").Value = 20
$ws.Range("P6").Value = 20

# --- Row 7 (full_task_roboto) -----------------------------------------------
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 1
$ws.Range("L7").Value = $true
$ws.Range("M7").Value = 5
$ws.Range("N7").Value = 10
$ws.Range("O7").Value = 12
$ws.Range("P7").Value = 12

# --- Row 8 (full_task_neuefrutigerworld) ------------------------------------
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 1
$ws.Range("L8").Value = $true
$ws.Range("M8").Value = 5
$ws.Range("N8").Value = 10
$ws.Range("O8").Value = 12
$ws.Range("P8").Value = 12

# --- Row 9 (full_task_eurostile) --------------------------------------------
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = 1
$ws.Range("L9").Value = $true
$ws.Range("M9").Value = 5
$ws.Range("N9").Value = 10
$ws.Range("O9").Value = 12
$ws.Range("P9").Value = 12

# --- Selection / view state --------------------------------------------------
[void]$ws.Range("F16").Select()
